$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$row = 33
$ws.Cells.Item($row, 2).Value = "SingleUseId36"
$ws.Cells.Item($row, 3).Value = "Iceland_45"
$ws.Cells.Item($row, 4).Value = "Center"
$ws.Cells.Item($row, 5).Value = "LTR"
$ws.Cells.Item($row, 6).Value = "MUSIC"
